# This workbook's update permutes the data rows (rows 2-51): the content
# that ends up in each row is the *old* content that used to live in a
# different row (same columns, A:AY). A few rows keep their own content
# (e.g. row 49). We therefore snapshot every row's old values first, then
# write the permuted set back in one shot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (1-based worksheet row numbers)
# i.e. new content of row $r == old content of row $map[$r]
$map = @{
  2=22;  3=8;   4=29;  5=13;  6=9;   7=2;   8=30;  9=15;  10=32;
  11=40; 12=20; 13=33; 14=35; 15=14; 16=6;  17=27; 18=25; 19=21;
  20=38; 21=11; 22=18; 23=34; 24=7;  25=36; 26=43; 27=17; 28=41;
  29=23; 30=16; 31=28; 32=31; 33=42; 34=24; 35=10; 36=5;  37=39;
  38=3;  39=44; 40=12; 41=4;  42=37; 43=19; 44=26; 45=48; 46=45;
  47=51; 48=50; 49=49; 50=47; 51=46
}

$firstRow = 2
$lastRow = 51
$nRows = $lastRow - $firstRow + 1
$firstCol = 1   # A
$lastCol = 51   # AY
$nCols = $lastCol - $firstCol + 1

$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))

# Snapshot all old values before we overwrite anything.
$old = $srcRange.Value()

# Some columns hold digit-only values that are nevertheless stored as
# text in the workbook: I (9, "Antal") has plain numbers typed as text,
# and Y (25)/AA (27) hold dates written as plain text (e.g.
# "2023-09-15"). If we assign such strings straight back through
# .Value, Excel "helpfully" reinterprets them as real numbers/dates.
# Force those destination columns to Text format first so the
# round-tripped value stays a literal string; we restore the default
# style afterwards so no stray formatting is left behind.
$forceTextCols = @(9, 25, 27)
foreach ($c in $forceTextCols) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $c), $ws.Cells.Item($lastRow, $c))
    $colRange.NumberFormat = "@"
}

# Build the permuted array to write back.
$new = New-Object 'object[,]' $nRows, $nCols
for ($r = 1; $r -le $nRows; $r++) {
    $destRow = $firstRow + $r - 1
    $srcRow = $map[$destRow]
    $srcIdx = $srcRow - $firstRow + 1
    for ($c = 1; $c -le $nCols; $c++) {
        $new[$r - 1, $c - 1] = $old[$srcIdx, $c]
    }
}

$srcRange.Value = $new

foreach ($c in $forceTextCols) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $c), $ws.Cells.Item($lastRow, $c))
    $colRange.Style = "Normal"
}
